# Rename naif_id as reg_id. Add confusing id check.
# Update based on decisions in 2024-01-24 meeting.
#
# The sheet "基本資料" gets a new column (J) added for a second kind of
# validation message, several rows of data are replaced/extended (rows 1-4),
# the stray artifact row (179) is removed, and a handful of cells get a
# gray "confusing id" highlight fill reusing the workbook's existing
# highlight style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基本資料")

# Drop the stray far-away artifact row entirely (it only held a styled,
# empty placeholder cell) before writing the new data.
$ws.Rows.Item(179).Delete()

# Clear out the old row contents/styles for the region we are rewriting.
$ws.Range("A1:J4").Clear()

$grayFill = 14540253        # 0x00DDDDDD -- existing "confusing id" highlight fill
$dateFmt  = "yyyy-mm-dd h:mm:ss"   # existing numFmtId 164

# ---- Row 1 ----
$ws.Cells.Item(1,1).Value = "Y"
$ws.Cells.Item(1,2).Value = 185403
$ws.Cells.Item(1,2).Interior.Color = $grayFill
$ws.Cells.Item(1,3).Value = "有問題"
$ws.Cells.Item(1,4).Value = 44609
$ws.Cells.Item(1,4).NumberFormat = $dateFmt
$ws.Cells.Item(1,5).Value = "Y121005地"
$ws.Cells.Item(1,6).Value = "Y146101趙"
$ws.Cells.Item(1,7).Value = 244310
$ws.Cells.Item(1,8).Value = "地趙"
$ws.Cells.Item(1,9).Value = "母"
$ws.Cells.Item(1,10).Value = "['不允許有相近耳號']"

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value = "L"
$ws.Cells.Item(2,2).Value = 188003
$ws.Cells.Item(2,2).Interior.Color = $grayFill
$ws.Cells.Item(2,3).Value = "L?"
$ws.Cells.Item(2,4).Value = 44665
$ws.Cells.Item(2,4).NumberFormat = $dateFmt
$ws.Cells.Item(2,7).Value = "無登"
$ws.Cells.Item(2,9).Value = "母"
$ws.Cells.Item(2,10).Value = "['不允許有相近耳號']"

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value = "L"
$ws.Cells.Item(3,2).Value = 190202
$ws.Cells.Item(3,4).Value = 44705
$ws.Cells.Item(3,4).NumberFormat = $dateFmt
$ws.Cells.Item(3,5).Value = "Y155009合"
$ws.Cells.Item(3,6).Value = "Y126104地"
$ws.Cells.Item(3,7).Value = "????"
$ws.Cells.Item(3,7).Interior.Color = $grayFill
$ws.Cells.Item(3,8).Value = "????"
$ws.Cells.Item(3,9).Value = "母"
$ws.Cells.Item(3,10).Value = "['登錄號不能含有非數字字元 ']"

# ---- Row 4 (new row) ----
$ws.Cells.Item(4,1).Value = "Y"
$ws.Cells.Item(4,2).Value = 137105
$ws.Cells.Item(4,2).Interior.Color = $grayFill
$ws.Cells.Item(4,3).Value = 137108
$ws.Cells.Item(4,4).NumberFormat = $dateFmt
$ws.Cells.Item(4,5).Value = "Y195207王"
$ws.Cells.Item(4,6).Value = "Y182001趙"
$ws.Cells.Item(4,7).Value = 238789
$ws.Cells.Item(4,8).Value = "王趙"
$ws.Cells.Item(4,9).Value = "母"
$ws.Cells.Item(4,10).Value = "['需要有子代的生日才能設定親代', '需要有子代的生日才能設定親代', '不允許有相近耳號']"
